$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws2 = $wb.Worksheets.Item("ARM")
$ws3 = $wb.Worksheets.Item("BSM")
$ws4 = $wb.Worksheets.Item("CRP")
$ws5 = $wb.Worksheets.Item("CUL")
$ws6 = $wb.Worksheets.Item("GSM")
$ws7 = $wb.Worksheets.Item("LTW")
$ws8 = $wb.Worksheets.Item("WVR")

# --- ALC ---
# row 62
$ws1.Range("H62").Value = 7489.886
$ws1.Range("I62").Value = 7103.3335
$ws1.Range("J62").Value = 9229.375
$ws1.Range("K62").Value = 7103.3335
$ws1.Range("L62").Value = 9229.375
$ws1.Range("M62").Value = -6479.3335
$ws1.Range("N62").Value = -10477.375
# row 65
$ws1.Range("H65").Value = 7489.886
$ws1.Range("I65").Value = 7103.3335
$ws1.Range("J65").Value = 9229.375
$ws1.Range("K65").Value = 35516.6675
$ws1.Range("L65").Value = 46146.875
$ws1.Range("M65").Value = -32396.6675
$ws1.Range("N65").Value = -52386.875
# row 113
$ws1.Range("H113").Value = 6799
$ws1.Range("I113").Value = 5000
$ws1.Range("J113").Value = 7998.3335
$ws1.Range("K113").Value = 5000
$ws1.Range("L113").Value = 7998.3335
$ws1.Range("M113").Value = -1746
$ws1.Range("N113").Value = -14506.3335
# row 135
$ws1.Range("H135").Value = 11006.597
$ws1.Range("I135").Value = 716.4583
$ws1.Range("J135").Value = 65887.336
$ws1.Range("K135").Value = 6448.1247
$ws1.Range("L135").Value = 592986.024
$ws1.Range("M135").Value = -3913.1247
$ws1.Range("N135").Value = -598056.024
# row 137
$ws1.Range("H137").Value = 14458.308
$ws1.Range("I137").Value = 22308
$ws1.Range("K137").Value = 66924
$ws1.Range("M137").Value = -64374
# row 138
$ws1.Range("H138").Value = 1878.1852
$ws1.Range("I138").Value = 1164.2727
$ws1.Range("J138").Value = 3388.3845
$ws1.Range("K138").Value = 3492.8181
$ws1.Range("L138").Value = 10165.1535
$ws1.Range("M138").Value = 1647.1819
$ws1.Range("N138").Value = -20445.1535

# --- ARM ---
# row 61
$ws2.Range("H61").Value = 6742.407
$ws2.Range("I61").Value = 1719.8334
$ws2.Range("J61").Value = 16787.555
$ws2.Range("K61").Value = 1719.8334
$ws2.Range("L61").Value = 16787.555
$ws2.Range("M61").Value = -1507.8334
$ws2.Range("N61").Value = -17211.555
# row 132
$ws2.Range("H132").Value = 2272.8542
$ws2.Range("I132").Value = 2077.6572
$ws2.Range("J132").Value = 2798.3845
$ws2.Range("K132").Value = 6232.971600000001
$ws2.Range("L132").Value = 8395.1535
$ws2.Range("M132").Value = -3702.971600000001
$ws2.Range("N132").Value = -13455.1535
# row 136
$ws2.Range("H136").Value = 6742.407
$ws2.Range("I136").Value = 1719.8334
$ws2.Range("J136").Value = 16787.555
$ws2.Range("K136").Value = 5159.5002
$ws2.Range("L136").Value = 50362.665
$ws2.Range("M136").Value = -2609.5002
$ws2.Range("N136").Value = -55462.665

# --- BSM ---
# row 99
$ws3.Range("H99").Value = 2042.2858
$ws3.Range("I99").Value = 1781.1818
$ws3.Range("J99").Value = 2999.6667
$ws3.Range("K99").Value = 1781.1818
$ws3.Range("L99").Value = 2999.6667
$ws3.Range("M99").Value = -283.1818000000001
$ws3.Range("N99").Value = -5995.6667

# --- CRP ---
# row 31
$ws4.Range("H31").Value = 5557641
$ws4.Range("I31").Value = 6251720.5
$ws4.Range("K31").Value = 6251720.5
$ws4.Range("M31").Value = -6251425.5
# row 34
$ws4.Range("H34").Value = 5557641
$ws4.Range("I34").Value = 6251720.5
$ws4.Range("K34").Value = 6251720.5
$ws4.Range("M34").Value = -6251518.5
# row 92
$ws4.Range("H92").Value = 0
$ws4.Range("J92").Value = 0
$ws4.Range("L92").Value = 0
$ws4.Range("N92").ClearContents()
# row 95
$ws4.Range("H95").Value = 47333.332
$ws4.Range("J95").Value = 51000
$ws4.Range("L95").Value = 51000
$ws4.Range("N95").Value = -56492
# row 96
$ws4.Range("H96").Value = 17999
$ws4.Range("J96").Value = 13598.8
$ws4.Range("L96").Value = 13598.8
$ws4.Range("N96").Value = -19090.8
# row 99
$ws4.Range("H99").Value = 5740.294
$ws4.Range("I99").Value = 5407.846
$ws4.Range("K99").Value = 5407.846
$ws4.Range("M99").Value = -3909.846
# row 126
$ws4.Range("H126").Value = 5740.294
$ws4.Range("I126").Value = 5407.846
$ws4.Range("K126").Value = 16223.538
$ws4.Range("M126").Value = -13753.538
# row 134
$ws4.Range("H134").Value = 1595.3334
$ws4.Range("I134").Value = 1453.8823
$ws4.Range("K134").Value = 4361.6469
$ws4.Range("M134").Value = -1826.6469

# --- CUL ---
# row 46
$ws5.Range("H46").Value = 899.6667
$ws5.Range("J46").Value = 899.6667
$ws5.Range("L46").Value = 2699.0001
$ws5.Range("N46").Value = -2881.0001

# --- GSM ---
# row 97
$ws6.Range("H97").Value = 3106.2856
$ws6.Range("I97").Value = 2578.5715
$ws6.Range("J97").Value = 3634
$ws6.Range("K97").Value = 2578.5715
$ws6.Range("L97").Value = 3634
$ws6.Range("M97").Value = -2082.5715
$ws6.Range("N97").Value = -4626
# row 132
$ws6.Range("H132").Value = 2124.3333
$ws6.Range("I132").Value = 2058.625
$ws6.Range("J132").Value = 2650
$ws6.Range("K132").Value = 6175.875
$ws6.Range("L132").Value = 7950
$ws6.Range("M132").Value = -3645.875
$ws6.Range("N132").Value = -13010

# --- LTW ---
# row 55
$ws7.Range("H55").Value = 1003.96295
$ws7.Range("I55").Value = 97.46154
$ws7.Range("J55").Value = 1845.7142
$ws7.Range("K55").Value = 97.46154
$ws7.Range("L55").Value = 1845.7142
$ws7.Range("M55").Value = 75.53846
$ws7.Range("N55").Value = -2191.7142
# row 100
$ws7.Range("H100").Value = 2510.7778
$ws7.Range("I100").Value = 2866.3333
$ws7.Range("J100").Value = 2333
$ws7.Range("K100").Value = 2866.3333
$ws7.Range("L100").Value = 2333
$ws7.Range("M100").Value = -2325.3333
$ws7.Range("N100").Value = -3415
# row 109
$ws7.Range("H109").Value = 38797.8
$ws7.Range("J109").Value = 38797.8
$ws7.Range("L109").Value = 38797.8
$ws7.Range("N109").Value = -41571.8
# row 132
$ws7.Range("H132").Value = 2211.7307
$ws7.Range("I132").Value = 2323.6597
$ws7.Range("J132").Value = 1159.6
$ws7.Range("K132").Value = 6970.9791
$ws7.Range("L132").Value = 3478.8
$ws7.Range("M132").Value = -4440.9791
$ws7.Range("N132").Value = -8538.799999999999
# row 136
$ws7.Range("H136").Value = 2122.926
$ws7.Range("I136").Value = 2059.898
$ws7.Range("J136").Value = 2740.6
$ws7.Range("K136").Value = 6179.694
$ws7.Range("L136").Value = 8221.799999999999
$ws7.Range("M136").Value = -3629.694
$ws7.Range("N136").Value = -13321.8

# --- WVR ---
# row 132
$ws8.Range("H132").Value = 4824229
$ws8.Range("I132").Value = 5336838.5
$ws8.Range("J132").Value = 5701
$ws8.Range("K132").Value = 16010515.5
$ws8.Range("L132").Value = 17103
$ws8.Range("M132").Value = -16007985.5
$ws8.Range("N132").Value = -22163
# row 136
$ws8.Range("H136").Value = 15447.117
$ws8.Range("I136").Value = 17814.023
$ws8.Range("J136").Value = 4401.5557
$ws8.Range("K136").Value = 53442.069
$ws8.Range("L136").Value = 13204.6671
$ws8.Range("M136").Value = -50892.069
$ws8.Range("N136").Value = -18304.6671
